$d = $word.ActiveDocument

$replacements = @(
    @{old="301×3=903";   new="241×5=1205"},
    @{old="505×4=2020";  new="363×3=1089"},
    @{old="459×5=2295";  new="867×4=3468"},
    @{old="502×7=3514";  new="795×4=3180"},
    @{old="695×3=2085";  new="326×5=1630"},
    @{old="732×2=1464";  new="672×5=3360"},
    @{old="410×8=3280";  new="665×7=4655"},
    @{old="591×9=5319";  new="849×3=2547"},
    @{old="696×2=1392";  new="755×3=2265"},
    @{old="206×8=1648";  new="910×5=4550"},
    @{old="771×7=5397";  new="684×3=2052"},
    @{old="949×3=2847";  new="238×2=476"},
    @{old="213×5=1065";  new="867×3=2601"},
    @{old="534×4=2136";  new="940×6=5640"},
    @{old="506×6=3036";  new="388×4=1552"},
    @{old="458×4=1832";  new="404×4=1616"},
    @{old="436×5=2180";  new="120×3=360"},
    @{old="702×8=5616";  new="818×4=3272"},
    @{old="602×7=4214";  new="505×2=1010"},
    @{old="960×7=6720";  new="688×6=4128"},
    @{old="701×3=2103";  new="778×4=3112"},
    @{old="764×2=1528";  new="526×8=4208"},
    @{old="604×3=1812";  new="230×8=1840"},
    @{old="292×2=584";   new="816×6=4896"},
    @{old="195×4=780";   new="879×9=7911"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.new, 2)
}
